# Update the Training Dashboard with the new progress snapshot dated 04-Nov-2025.
# For each data row (3 through 30):
#   - Column H (PERIOD TO EXPIRE) is decremented by 1 day.
#   - Column I (LAST UPDATE) is changed from 03-Nov-2025 to 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 30; $r++) {
    $periodCell = $ws.Cells.Item($r, 8)   # column H: PERIOD TO EXPIRE
    $updateCell = $ws.Cells.Item($r, 9)   # column I: LAST UPDATE

    $currentPeriod = $periodCell.Value2
    if ($currentPeriod -ne $null -and $currentPeriod -ne "") {
        $periodCell.Value = $currentPeriod - 1
    }

    $currentUpdate = $updateCell.Value2
    if ($currentUpdate -eq "03-Nov-2025") {
        # Assign via a text formula and paste the computed value back so Excel
        # keeps storing this as plain text (matching the workbook's existing
        # convention for date-looking cells) instead of auto-converting the
        # string into a date serial number, which would also disturb the
        # cell's number format/style.
        $updateCell.Formula = '="04-Nov-2025"'
        $updateCell.Copy() | Out-Null
        $updateCell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    }
}

$excel.CutCopyMode = $false
